$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Content edits, issued in the same order the author likely made them so
# that newly-created shared-string entries land at the same indices as the
# target workbook. ---

# 1) Sheet1 ("Commands"): add a note about further Alter-table topics to explore
$ws1.Range("D13").Value2 = "MYSQL: More topics to explore on Alter table.`n# Alter table name`n# Alter column name  `n# Alter table - add column`n# Alter table - delete column"

# 2) Sheet2 ("Topics to Learn"): mark "establish connection" / "creating a db" topics done
$ws2.Range("C2").Value2 = "Done"
$ws2.Range("D2").Value2 = "Done"
$ws2.Range("C4").Value2 = "Done"
$ws2.Range("D4").Value2 = "Done"

# 3) Update the "data types" row with the real datatype categories
$ws2.Range("B10").Value2 = "Numeric, String, Data time"

# 4) New "Alter table" topic row with its practice snippet
$ws2.Range("A11").Value2 = "ALTER Table "
$ws2.Range("B11").Value2 = "# Alter table name`n# Alter column name  `n# Alter table - add column`n# Alter table - delete column"
$ws2.Range("B11").WrapText = $true

# 5) New "mycursor" methods rows
$ws2.Range("B12").Value2 = "mycursor.rowcount"
$ws2.Range("B13").Value2 = "mycursor.fetchone"
$ws2.Range("A12").Value2 = "Explore mycursor methods"

# 6) Sheet1: replace the old "insert many rows" example with a practical single insert example
$ws1.Range("B15").Value2 = "sql = `"INSERT INTO friends (name, email, mobile) VALUES (%s,%s,%s)`"`nval = ('Aditya', 'aditya@gmail.com', '99999')`nmycursor.execute(sql,val)`nmydb.commit()`nmycursor.rowcount`nprint(`"1 record inserted, ID:`", mycursor.lastrowid)"

# --- Row-height touch-ups to mirror Excel's own auto-fit results for the
# edited rows ---
$ws1.Rows(13).RowHeight = 75
$ws1.Rows(15).RowHeight = 105
$ws2.Rows(11).RowHeight = 60

# --- View/selection updates ---
$ws2.Range("B14").Select()
$ws1.Activate()
$ws1.Range("B15").Select()
